$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data driving the "10 inch setup" row (row 7)
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 5.4499999999999282
$ws.Range("I7").Value = 4.2199999999999545

# Update the active selection to I7
$ws.Range("I7").Select()
